$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.483.09"
$ws.Range("E2").Value = "  +2.18%  "

$ws.Range("D3").Value = "1.864.37"
$ws.Range("E3").Value = "  +2.61%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.45%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.45%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4669"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.96%  "

$ws.Range("E8").Value = "  +2.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07378"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8895"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.73%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07967"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.00"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.42%  "

$ws.Range("D13").Value = "1.871.68"
$ws.Range("E13").Value = "  +3.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.428"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.95%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.600"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.70%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "92.69"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.005"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.37%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008954"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.93%  "

$ws.Range("E19").Value = "  -0.39%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.37%  "

$ws.Range("D21").Value = "27.496.52"
$ws.Range("E21").Value = "  +2.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.169"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.74%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.61%  "

$ws.Range("D24").Value = "2.090.27"
$ws.Range("E24").Value = "  +3.50%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.91%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.874"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.087"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.80%  "

$ws.Range("E29").Value = "  +1.56%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08913"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.52%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7535"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.64%  "

$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.008"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.36%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.162"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.35%  "

$ws.Range("E35").Value = "  +1.87%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.643"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01970"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.77%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.082"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05287"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.82%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.985"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.41%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.173"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.50%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5219"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.64%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1646"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.50%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.365"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.71%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4886"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.88%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.86%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.003"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.51%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.51%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.661"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.93%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06259"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.06%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.56%  "
